# fix(stats): add stats to export
#
# Inserts three new "reexpedie" stat rows into the "Stats structure" sheet:
#   - "Colis rééxpédiés"            (before "Courriers enregistrés")
#   - "Courriers réexpédiés"        (after  "Courriers remis")
#   - "Avis de passage réexpédiés"  (after  "Avis de passage remis")
#
# Original rows 140..149 (1-based, col B text / col C blank):
#   140 Courriers enregistrés
#   141 Courriers remis
#   142 Avis de passage enregistrés
#   143 Avis de passage remis
#   144 Passages enregistrés
#   145 Passages avec remise de courrier
#   146 Passages sans remise de courrier
#   147 Connexion au portail usager
#   148 (blank)
#   149 (blank)
#
# Target rows 140..152:
#   140 Colis rééxpédiés            <- NEW
#   141 Courriers enregistrés
#   142 Courriers remis
#   143 Courriers réexpédiés        <- NEW
#   144 Avis de passage enregistrés
#   145 Avis de passage remis
#   146 Avis de passage réexpédiés  <- NEW
#   147 Passages enregistrés
#   148 Passages avec remise de courrier
#   149 Passages sans remise de courrier
#   150 Connexion au portail usager
#   151 (blank)
#   152 (blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert the three blank rows bottom-up (so the row numbers used for the
# not-yet-processed insertions above stay valid while we work).
$ws.Rows.Item(144).Insert()
$ws.Rows.Item(142).Insert()
$ws.Rows.Item(140).Insert()

# Now fill in the new cells' text, in top-to-bottom (final sheet string)
# order so the new shared-string entries get appended in the same order
# as the target workbook (114=Colis, 115=Courriers, 116=Avis de passage).
$ws.Cells.Item(140, 2).Value = "Colis rééxpédiés"
$ws.Cells.Item(143, 2).Value = "Courriers réexpédiés"
$ws.Cells.Item(146, 2).Value = "Avis de passage réexpédiés"

# Match the row styling used by all the other stat rows in this block
# (left-aligned label style "16" on column B, style "6" on column C -- both
# already carried over onto the freshly inserted rows by Insert(), but make
# sure the row height bucket used across this section is consistent).
$ws.Rows.Item(140).RowHeight = 16
$ws.Rows.Item(143).RowHeight = 16
$ws.Rows.Item(146).RowHeight = 16

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Activate()
$ws.Range("E129").Select()
